# RPA datasets push 2024-06-22
# Insert a new data row at row 2 (new IPO/SPAC entry: KB제29호스팩, listed 2024-06-21),
# pushing the existing rows 2-24 down to rows 3-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (the first data row, right below the header).
$ws.Rows.Item(2).Insert()

# The inserted row copies the header row's bold/border formatting by default;
# strip that so the new row matches the unstyled data rows beneath it.
$ws.Rows.Item(2).ClearFormats()

# Fill in the new record's values. Plain text/number columns can be assigned
# directly.
$ws.Cells.Item(2, 2).Value = "KB제29호스팩"
$ws.Cells.Item(2, 3).Value = "코스닥"
$ws.Cells.Item(2, 4).Value = 120
$ws.Cells.Item(2, 5).Value = "KB"
$ws.Cells.Item(2, 6).Value = 120
$ws.Cells.Item(2, 7).Value = "-"
$ws.Cells.Item(2, 8).Value = "-"
$ws.Cells.Item(2, 9).Value = "-"
$ws.Cells.Item(2, 10).Value = "-"
$ws.Cells.Item(2, 11).Value = "대표"
$ws.Cells.Item(2, 12).Value = "-"
$ws.Cells.Item(2, 13).Value = 2000
$ws.Cells.Item(2, 14).Value = 100
$ws.Cells.Item(2, 17).Value = 4500000

# Date-like text ("YYYY-MM-DD") would otherwise auto-convert to a real date
# serial when assigned directly; a leading apostrophe forces literal text
# entry (as typing it in the Excel UI would), matching how the rest of the
# sheet's date columns (A, O, P) store their values as plain text.
$ws.Cells.Item(2, 1).Value = "'2024-06-21"
$ws.Cells.Item(2, 15).Value = "'2024-06-11"
$ws.Cells.Item(2, 16).Value = "'2024-06-14"

# The quote-prefix entry marks the cells with a "stored as text" style;
# clear formatting once more so the row ends up with no style overrides,
# consistent with every other data row in the sheet.
$ws.Rows.Item(2).ClearFormats()
